$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting nothing below it (sheet only has
# rows 1-3 so far), which extends the used range down to row 4 and keeps
# row 3's formatting available to copy from if needed.
$ws.Rows.Item(4).Insert()

# --- Row 4: brand new row, identical to the ORIGINAL row 2 data ---
$ws.Range("A4").Value2 = 1
$ws.Range("B4").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value2 = "Arica y Parinacota"
$ws.Range("D4").Value2 = 44167
$ws.Range("D4").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E4").Value2 = 15
$ws.Range("F4").Value2 = "Fruta"
$ws.Range("G4").Value2 = 100107
$ws.Range("H4").Value2 = "Otros"
$ws.Range("I4").Value2 = 100107002
$ws.Range("J4").Value2 = "Chirimoya"
$ws.Range("K4").Value2 = "Cultivar IV Región"
$ws.Range("L4").Value2 = "Segunda"
$ws.Range("M4").Value2 = 200
$ws.Range("N4").Value2 = 18000
$ws.Range("O4").Value2 = 19000
$ws.Range("P4").Value2 = 18500
$ws.Range("Q4").Value2 = "$/caja 13 kilos"
$ws.Range("R4").Value2 = "Región de Coquimbo"
$ws.Range("S4").Value2 = 1423
$ws.Range("T4").Value2 = 13

# --- Row 2: update date + price fields to match the figures that used to
# live in row 3 (Segunda quality, same volume) ---
$ws.Range("D2").Value2 = 44160
$ws.Range("N2").Value2 = 19000
$ws.Range("O2").Value2 = 20000
$ws.Range("P2").Value2 = 19500
$ws.Range("S2").Value2 = 1500

# --- Row 3: now holds a brand new "Primera" quality entry ---
$ws.Range("D3").Value2 = 44441
$ws.Range("L3").Value2 = "Primera"
$ws.Range("M3").Value2 = 100
$ws.Range("N3").Value2 = 29000
$ws.Range("O3").Value2 = 30000
$ws.Range("P3").Value2 = 29500
$ws.Range("Q3").Value2 = "$/caja 12 kilos"
$ws.Range("S3").Value2 = 2458
$ws.Range("T3").Value2 = 12
